$d = $word.ActiveDocument
$rng = $d.Content
$null = $rng.Find.Execute("- Introdução à Ciência e Tecnologia de Polímeros.", $true, $false, $false, $false, $false, $true, 1, $false, "- Introdução à Ciência e Tecnologia de Polímeros.^l", 2)
$null = $rng.Find.Execute("- Considerações sobre a síntese de materiais poliméricos avançados.", $true, $false, $false, $false, $false, $true, 1, $false, "- Considerações sobre a síntese de materiais poliméricos avançados.^l", 2)
$null = $rng.Find.Execute("- Polímeros com propriedades mecânicas excepcionais.", $true, $false, $false, $false, $false, $true, 1, $false, "- Polímeros com propriedades mecânicas excepcionais.^l", 2)
$null = $rng.Find.Execute("- Polímeros com memória de forma.", $true, $false, $false, $false, $false, $true, 1, $false, "- Polímeros com memória de forma.^l", 2)
$null = $rng.Find.Execute("- Aplicações selecionadas de polímeros em:", $true, $false, $false, $false, $false, $true, 1, $false, "- Aplicações selecionadas de polímeros em:^l", 2)
$null = $rng.Find.Execute("- transporte de fármacos (drug delivery);", $true, $false, $false, $false, $false, $true, 1, $false, "- transporte de fármacos (drug delivery);^l", 2)
$null = $rng.Find.Execute("- transformação de células (transfecção);", $true, $false, $false, $false, $false, $true, 1, $false, "- transformação de células (transfecção);^l", 2)
$null = $rng.Find.Execute("- próteses de base polimérica para uso em humanos;", $true, $false, $false, $false, $false, $true, 1, $false, "- próteses de base polimérica para uso em humanos;^l", 2)
$null = $rng.Find.Execute("- nanorreatores/catálise;", $true, $false, $false, $false, $false, $true, 1, $false, "- nanorreatores/catálise;^l", 2)
$null = $rng.Find.Execute("- descontaminação de corpos d ́água e outras aplicações relacionadas ao meio-ambiente;", $true, $false, $false, $false, $false, $true, 1, $false, "- descontaminação de corpos d ́água e outras aplicações relacionadas ao meio-ambiente;^l", 2)
$null = $rng.Find.Execute("- eletrônica/polímeros condutores;", $true, $false, $false, $false, $false, $true, 1, $false, "- eletrônica/polímeros condutores;^l", 2)
$null = $rng.Find.Execute("- agricultura", $true, $false, $false, $false, $false, $true, 1, $false, "- agricultura^l", 2)
$null = $rng.Find.Execute("- revestimentos (coatings) ativos de superfícies.", $true, $false, $false, $false, $false, $true, 1, $false, "- revestimentos (coatings) ativos de superfícies.^l", 2)
$null = $rng.Find.Execute("- recuperação avançada de petróleo.", $true, $false, $false, $false, $false, $true, 1, $false, "- recuperação avançada de petróleo.^l", 2)
$null = $rng.Find.Execute("- Polímeros foto/bio/oxidegradáveis", $true, $false, $false, $false, $false, $true, 1, $false, "- Polímeros foto/bio/oxidegradáveis^l", 2)
$null = $rng.Find.Execute("- Polímeros e os seres vivos/Biopolímeros.", $true, $false, $false, $false, $false, $true, 1, $false, "- Polímeros e os seres vivos/Biopolímeros.^l", 2)
